$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 16773.143
$ws.Range("J53").Value = 603
$ws.Range("L53").Value = 603
$ws.Range("N53").Value = -1877
$ws.Range("H70").Value = 7774.4375
$ws.Range("I70").Value = 5663.2856
$ws.Range("K70").Value = 16989.8568
$ws.Range("M70").Value = -16719.8568
$ws.Range("H73").Value = 7774.4375
$ws.Range("I73").Value = 5663.2856
$ws.Range("K73").Value = 16989.8568
$ws.Range("M73").Value = -16053.8568
$ws.Range("H107").Value = 765.375
$ws.Range("I107").Value = 571.8889
$ws.Range("K107").Value = 571.8889
$ws.Range("M107").Value = 1348.1111
$ws.Range("H132").Value = 1143.25
$ws.Range("H135").Value = 820.3333
$ws.Range("I135").Value = 771.6667
$ws.Range("K135").Value = 6945.0003
$ws.Range("M135").Value = -4410.0003
$ws.Range("J141").Value = 3821
$ws.Range("L141").Value = 11463
$ws.Range("N141").Value = -21823

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1662495.9
$ws.Range("I2").Value = 2326804.2
$ws.Range("K2").Value = 2326804.2
$ws.Range("M2").Value = -2326691.2
$ws.Range("H32").Value = 2560.7349
$ws.Range("I32").Value = 1977.9546
$ws.Range("J32").Value = 4823.294
$ws.Range("K32").Value = 1977.9546
$ws.Range("L32").Value = 4823.294
$ws.Range("M32").Value = -1690.9546
$ws.Range("N32").Value = -5397.294
$ws.Range("H37").Value = 13950
$ws.Range("J37").Value = 13950
$ws.Range("L37").Value = 13950
$ws.Range("N37").Value = -14496
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = $null
$ws.Range("H61").Value = 3480.7778
$ws.Range("I61").Value = 1099.75
$ws.Range("K61").Value = 1099.75
$ws.Range("M61").Value = -887.75
$ws.Range("H74").Value = 1158.826
$ws.Range("I74").Value = 737.4
$ws.Range("K74").Value = 737.4
$ws.Range("M74").Value = 136.6
$ws.Range("H77").Value = 1158.826
$ws.Range("I77").Value = 737.4
$ws.Range("K77").Value = 3687
$ws.Range("M77").Value = 681
$ws.Range("H93").Value = 29800
$ws.Range("J93").Value = 29800
$ws.Range("L93").Value = 29800
$ws.Range("N93").Value = -34792
$ws.Range("H110").Value = 1110.6
$ws.Range("I110").Value = 888.25
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 888.25
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1156.75
$ws.Range("N110").Value = -6090
$ws.Range("H116").Value = 1662495.9
$ws.Range("I116").Value = 2326804.2
$ws.Range("K116").Value = 2326804.2
$ws.Range("M116").Value = -2324510.2
$ws.Range("H136").Value = 3480.7778
$ws.Range("I136").Value = 1099.75
$ws.Range("K136").Value = 3299.25
$ws.Range("M136").Value = -749.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1662495.9
$ws.Range("I3").Value = 2326804.2
$ws.Range("K3").Value = 2326804.2
$ws.Range("M3").Value = -2326690.2
$ws.Range("H20").Value = 1970.3334
$ws.Range("I20").Value = 1925.8462
$ws.Range("K20").Value = 1925.8462
$ws.Range("M20").Value = -1678.8462
$ws.Range("H25").Value = 55000
$ws.Range("J25").Value = 55000
$ws.Range("L25").Value = 55000
$ws.Range("N25").Value = -55470
$ws.Range("H94").Value = 1072.375
$ws.Range("I94").Value = 1082.7142
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1082.7142
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -631.7141999999999
$ws.Range("N94").Value = -1902
$ws.Range("H96").Value = 10891.667
$ws.Range("I96").Value = 10891.667
$ws.Range("K96").Value = 10891.667
$ws.Range("M96").Value = -8145.666999999999
$ws.Range("H105").Value = 2924.9167
$ws.Range("I105").Value = 2981.7273
$ws.Range("K105").Value = 2981.7273
$ws.Range("M105").Value = -1234.7273
$ws.Range("H134").Value = 5925.615
$ws.Range("I134").Value = 6118.892
$ws.Range("K134").Value = 18356.676
$ws.Range("M134").Value = -15821.676

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1742.683
$ws.Range("I31").Value = 901.0769
$ws.Range("J31").Value = 2133.4285
$ws.Range("K31").Value = 901.0769
$ws.Range("L31").Value = 2133.4285
$ws.Range("M31").Value = -606.0769
$ws.Range("N31").Value = -2723.4285
$ws.Range("H34").Value = 1742.683
$ws.Range("I34").Value = 901.0769
$ws.Range("J34").Value = 2133.4285
$ws.Range("K34").Value = 901.0769
$ws.Range("L34").Value = 2133.4285
$ws.Range("M34").Value = -699.0769
$ws.Range("N34").Value = -2537.4285
$ws.Range("H88").Value = 5000
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 5000
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 5000
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 5000
$ws.Range("N91").Value = -7808

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 527.5
$ws.Range("J2").Value = 487.69232
$ws.Range("L2").Value = 2926.15392
$ws.Range("N2").Value = -3152.15392
$ws.Range("H33").Value = 210.77777
$ws.Range("I33").Value = 112
$ws.Range("J33").Value = 260.16666
$ws.Range("K33").Value = 672
$ws.Range("L33").Value = 1560.99996
$ws.Range("M33").Value = -389
$ws.Range("N33").Value = -2126.99996
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = $null
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = 0
$ws.Range("H131").Value = 8076867
$ws.Range("I131").Value = 250000510
$ws.Range("K131").Value = 750001530
$ws.Range("M131").Value = -749996490

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3852.75
$ws.Range("I113").Value = 5905.5
$ws.Range("K113").Value = 5905.5
$ws.Range("M113").Value = -3735.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5905.4
$ws.Range("I7").Value = 3626
$ws.Range("J7").Value = 7425
$ws.Range("K7").Value = 3626
$ws.Range("L7").Value = 7425
$ws.Range("M7").Value = -3514
$ws.Range("N7").Value = -7649
$ws.Range("H55").Value = 225
$ws.Range("J55").Value = 171
$ws.Range("L55").Value = 171
$ws.Range("N55").Value = -517
$ws.Range("H93").Value = 1013
$ws.Range("I93").Value = 819.6
$ws.Range("J93").Value = 1980
$ws.Range("K93").Value = 819.6
$ws.Range("L93").Value = 1980
$ws.Range("M93").Value = 428.4
$ws.Range("N93").Value = -4476
$ws.Range("H126").Value = 5905.4
$ws.Range("I126").Value = 3626
$ws.Range("J126").Value = 7425
$ws.Range("K126").Value = 10878
$ws.Range("L126").Value = 22275
$ws.Range("M126").Value = -8408
$ws.Range("N126").Value = -27215
$ws.Range("H132").Value = 3019.8635
$ws.Range("I132").Value = 1033.5217
$ws.Range("K132").Value = 3100.5651
$ws.Range("M132").Value = -570.5650999999998
$ws.Range("H136").Value = 2976.4773
$ws.Range("I136").Value = 2083.0645
$ws.Range("J136").Value = 5106.923
$ws.Range("K136").Value = 6249.193499999999
$ws.Range("L136").Value = 15320.769
$ws.Range("M136").Value = -3699.193499999999
$ws.Range("N136").Value = -20420.769

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12079240
$ws.Range("J136").Value = 1574.7667
$ws.Range("L136").Value = 4724.300099999999
$ws.Range("N136").Value = -9824.3001
